$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows with changed data.
# D-column values are textual (often contain multiple "." as thousand separators),
# so they are written with a leading apostrophe to force text storage, matching
# the original inlineStr type instead of being auto-converted to a number.
$ws.Range("D2").Value = "'56.984.82"
$ws.Range("E2").Value = "  +0.84%  "

$ws.Range("D3").Value = "'2.342.46"
$ws.Range("E3").Value = "  -0.04%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'518.78"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("D6").Value = "'135.88"
$ws.Range("E6").Value = "  +2.02%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  +0.41%  "

$ws.Range("D9").Value = "'2.353.32"
$ws.Range("E9").Value = "  +0.38%  "

$ws.Range("E10").Value = "  -0.34%  "

$ws.Range("D11").Value = "'5.42"
$ws.Range("E11").Value = "  +5.12%  "

$ws.Range("E12").Value = "  -1.51%  "

$ws.Range("E13").Value = "  +0.08%  "

$ws.Range("E14").Value = "  -0.16%  "

$ws.Range("D15").Value = "'2.763.50"
$ws.Range("E15").Value = "  +1.04%  "

$ws.Range("D16").Value = "'56.972.88"
$ws.Range("E16").Value = "  +0.95%  "

$ws.Range("E17").Value = "  -0.18%  "

$ws.Range("D18").Value = "'2.358.01"
$ws.Range("E18").Value = "  -0.57%  "

$ws.Range("D19").Value = "'10.62"
$ws.Range("E19").Value = "  +0.64%  "

$ws.Range("D20").Value = "'327.17"
$ws.Range("E20").Value = "  +1.58%  "

$ws.Range("E21").Value = "  -0.34%  "

$ws.Range("D22").Value = "'6.75"
$ws.Range("E22").Value = "  +1.31%  "

$ws.Range("E23").Value = "  -0.20%  "

$ws.Range("D24").Value = "'61.16"
$ws.Range("E24").Value = "  +0.25%  "

$ws.Range("D25").Value = "'0.166"
$ws.Range("E25").Value = "  +4.97%  "

$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  +0.42%  "

$ws.Range("D27").Value = "'8.01"
$ws.Range("E27").Value = "  +4.72%  "

$ws.Range("E28").Value = "  +9.75%  "

$ws.Range("D29").Value = "'170.33"
$ws.Range("E29").Value = "  -1.02%  "

$ws.Range("E30").Value = "  +3.04%  "

$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("D33").Value = "'18.57"
$ws.Range("E33").Value = "  +1.06%  "

$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").Value = "'0.996"
$ws.Range("E35").Value = "  -0.39%  "

$ws.Range("E36").Value = "  +0.82%  "

$ws.Range("D37").Value = "'0.915"
$ws.Range("E37").Value = "  -1.20%  "

$ws.Range("E38").Value = "  +1.16%  "

$ws.Range("E39").Value = "  +3.64%  "

$ws.Range("D40").Value = "'38.40"
$ws.Range("E40").Value = "  +2.58%  "

$ws.Range("D41").Value = "'149.17"
$ws.Range("E41").Value = "  +7.66%  "

$ws.Range("D42").Value = "'0.384"
$ws.Range("E42").Value = "  +0.24%  "

$ws.Range("D43").Value = "'3.65"
$ws.Range("E43").Value = "  +0.84%  "

$ws.Range("D46").Value = "'0.0936"
$ws.Range("E46").Value = "  +1.17%  "

$ws.Range("E47").Value = "  -0.65%  "

$ws.Range("D48").Value = "'0.563"
$ws.Range("E48").Value = "  +1.37%  "

$ws.Range("E49").Value = "  +1.96%  "

$ws.Range("E50").Value = "  +5.83%  "

$ws.Range("D51").Value = "'0.382"
$ws.Range("E51").Value = "  -0.08%  "

# Rows 44 and 45 swap content: Bittensor moves up to row 44 (previously RenderToken),
# RenderToken moves down to row 45 (previously Bittensor), each with new price/volume data.
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "'279.91"
$ws.Range("E44").Value = "  +4.41%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'5.21"
$ws.Range("E45").Value = "  +2.66%  "
